$d = $word.ActiveDocument

# En dash (U+2013) is used elsewhere in the document for these "Key - Value"
# lines (e.g. "ART " + dash + " Art Assets").
$dash = [char]0x2013

# Find the paragraph that reads "DOCUMENTATION - Documentation Documents";
# the very next paragraph is the empty one that needs to become two new
# lines followed by an (empty) paragraph holding the relocated bookmark.
$docIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "DOCUMENTATION $dash Documentation Documents" -or $t.StartsWith("DOCUMENTATION")) {
        $docIndex = $i
        break
    }
}

$r = $d.Paragraphs.Item($docIndex + 1).Range
$r.InsertBefore("CONTROL $dash Animation Controller" + [char]13 + "PREFAB $dash A Unity Prefab" + [char]13)

# InsertBefore grows the original range to cover the newly-inserted text,
# so re-fetch the paragraph by (now shifted) index to get the still-empty
# trailing paragraph - that's where the bookmark needs to end up. Word
# only ever keeps a single _GoBack bookmark, so re-adding it relocates it
# and implicitly drops the old one that sat on the Unity Version line.
$target = $d.Paragraphs.Item($docIndex + 3)
$d.Bookmarks.Add("_GoBack", $target.Range)
